$d = $word.ActiveDocument

# --- Update the three section header labels to the new uppercase wording ---
$d.Content.Find.Execute("A. Formarea profesionala a beneficiarilor", $true, $false, $false, $false, $false,
                         $true, 1, $false, "A. PERFORMANTA IN FORMAREA PROFESIONALA A BENEFICIARILOR", 2)

$d.Content.Find.Execute("B. Performanta in cercetarea stiintifica", $true, $false, $false, $false, $false,
                         $true, 1, $false, "B. PERFORMANTA IN CERCETAREA STIINTIFICA", 2)

$d.Content.Find.Execute("C. Performanta privind participarea la dezvoltarea institutionala", $true, $false, $false, $false, $false,
                         $true, 1, $false, "C. PERFORMANTA PRIVIND PARTICIPAREA LA DEZVOLTAREA INSTITUTIONALA", 2)

# --- Update the performance figures for the first data row (First Administrator) ---
$t = $d.Tables.Item(1)
$row = $t.Rows.Item(3)
$row.Cells.Item(3).Range.Text = "91.22"
$row.Cells.Item(4).Range.Text = "A"
$row.Cells.Item(5).Range.Text = "16.33"
$row.Cells.Item(6).Range.Text = "B"
$row.Cells.Item(7).Range.Text = "19.0"
$row.Cells.Item(8).Range.Text = "A"
$row.Cells.Item(9).Range.Text = "3.67"

# --- Remove the now-obsolete second data row (First User) entirely ---
$t.Rows.Item(4).Delete()

Write-Output "done"
